$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.910.07"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.703.31"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.86%  "
$ws.Range("D5").Value = "315.72"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").Value = "0.4071"
$ws.Range("E7").Value = "  +3.40%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "1.004"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").Value = "53.79"
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("D11").Value = "1.467"
$ws.Range("E11").Value = "  -2.65%  "
$ws.Range("D12").Value = "0.08816"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "25.77"
$ws.Range("E13").Value = "  +4.95%  "
$ws.Range("D14").Value = "7.489"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "8.045"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "0.00001351"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "1.729.72"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("E18").Value = "  -2.81%  "
$ws.Range("D19").Value = "0.07170"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").Value = "20.95"
$ws.Range("E20").Value = "  +5.48%  "
$ws.Range("D21").Value = "7.233"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").Value = "14.58"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").Value = "24.901.99"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "2.326"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").Value = "6.746"
$ws.Range("E26").Value = "  +29.10%  "
$ws.Range("D27").Value = "2.881"
$ws.Range("E27").Value = "  -5.17%  "
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("D29").Value = "164.72"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "144.87"
$ws.Range("E30").Value = "  +4.93%  "
$ws.Range("D31").Value = "8.242"
$ws.Range("E31").Value = "  -4.44%  "
$ws.Range("D32").Value = "2.270"
$ws.Range("E32").Value = "  +14.12%  "
$ws.Range("D33").Value = "1.918.07"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("D34").Value = "0.08776"
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("E35").Value = "  +10.50%  "
$ws.Range("D36").Value = "7.302"
$ws.Range("E36").Value = "  -4.36%  "
$ws.Range("D37").Value = "1.018"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").Value = "0.2843"
$ws.Range("E38").Value = "  +3.79%  "
$ws.Range("D39").Value = "0.8505"
$ws.Range("E39").Value = "  +8.47%  "
$ws.Range("D40").Value = "10.98"
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("D41").Value = "0.09442"
$ws.Range("E41").Value = "  +3.50%  "
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").Value = "1.470"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "17.70"
$ws.Range("E44").Value = "  +6.56%  "
$ws.Range("D45").Value = "2.725"
$ws.Range("E45").Value = "  +5.08%  "
$ws.Range("D46").Value = "0.7433"
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "1.400"
$ws.Range("E48").Value = "  +5.47%  "
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").Value = "142.03"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("E51").Value = "  +4.26%  "